$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '275.15'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.47%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.41'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.41%'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.98%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06285'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.95%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.923'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.32%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.312'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '38.21%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8748'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.82%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1523'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.49%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05044'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.30%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07476'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.08%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02923'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-6.84%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09052'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.12%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001566'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.72%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006344'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.95%'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2.94%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.448'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.94%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.319'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.04%'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.02%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1319'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.74%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.930'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.78%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04400'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.86%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001171'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.28%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003826'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '5.65%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001200'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.08%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001940'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '14.62%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04105'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.46%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006905'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4.33%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1171'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.51%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001930'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-17.82%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01123'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-10.10%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005210'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.02%'

$ws.Range("B46").Value = 'BOLO'

$ws.Range("C46").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.490'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-37.38%'

$ws.Range("B47").Value = 'CoinbaseStockToken'

$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.02003'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-10.87%'
